$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename headers D1 and E1 (casing fix metadata4Ing -> metadata4ing)
$ws.Range("D1").Value = "metadata4ing_IRI"
$ws.Range("E1").Value = "metadata4ing_DESC"

# Add new header F1 (copy formatting from an existing header cell, then set text)
$ws.Range("C1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "metadata4ing_DEF"

# Fill new column F values for data rows
$ws.Range("F2").Value = "[]"
$ws.Range("F3").Value = "[]"
$ws.Range("F4").Value = "[]"
$ws.Range("F5").Value = "[]"
$ws.Range("F6").Value = "[locstr('A role is the function of an entity or agent with respect to an activity, in the context of a usage, generation, invalidation, association, start, and end.', 'en')]"
$ws.Range("F7").Value = "[locstr('Property, i.e., a variable that can be considered (as an investigated property, i.e., target property) within a processing step', 'en')]"
$ws.Range("F8").Value = "['To say that b is a realizable entity is to say that b is a specifically dependent continuant that inheres in some independent continuant which is not a spatial region and is of a type instances of which are realized in processes of a correlated type.´[BFO]', 'To say that b is a realizable entity is to say that b is a specifically dependent continuant that inheres in some independent continuant which is not a spatial region and is of a type instances of which are realized in processes of a correlated type. (axiom label in BFO2 Reference: [058-002])']"
